# Daily attendance processing - 2026-01-08 08:43:19
#
# Applies the day's processing pass over the "Session Analysis Results"
# sheet:
#  1. Bumps the summary "Missing Sessions" (L7) and "Pending Sessions" (L8)
#     counters.
#  2. Re-orders the "Recorded By" list for every session that was recorded
#     by both the System and dnasr281@gmail.com, putting "System" first.
#  3. Updates the per-group Pending/Missing session tallies (P21:Q26).
#  4. Rolls sessions dated 08/01/2026 (today) from the "Pending" (future)
#     look into "Not Recorded" (past-due, unrecorded) once their date has
#     elapsed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Summary counters (K7/L7 "Missing Sessions", K8/L8 "Pending Sessions") ---
$ws.Range("L7").Value = 33
$ws.Range("L8").Value = 72

# --- 2. "Recorded By" re-ordering: "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com" ---
$recordedByRows = @(8,9,10,12,14,15,17,18,34,35,36,38,40,41,43,44,60,61,62,64,66,67,69,70,86,87,88,90,92,93,95,96,112,113,114,116,118,119,121,122,138,139,140,142,144,145,147,148,164,167,170,174,191,194,197,201,218,221,224,228,245,248,251,255,272,275,278,282,299,302,305,309)
foreach ($row in $recordedByRows) {
    $cell = $ws.Cells.Item($row, 7)
    if ($cell.Value() -eq "dnasr281@gmail.com, System") {
        $cell.Value = "System, dnasr281@gmail.com"
    }
}

# --- 3. Per-group Pending/Missing tallies (columns P and Q, rows 21-26) ---
$pqUpdates = @{
    21 = @{ P = 4; Q = 6 }
    22 = @{ P = 4; Q = 6 }
    23 = @{ P = 4; Q = 6 }
    24 = @{ P = 5; Q = 6 }
    25 = @{ P = 4; Q = 6 }
    26 = @{ P = 4; Q = 6 }
}
foreach ($row in $pqUpdates.Keys) {
    $ws.Cells.Item($row, 16).Value = $pqUpdates[$row].P
    $ws.Cells.Item($row, 17).Value = $pqUpdates[$row].Q
}

# --- 4. Roll today's (08/01/2026) sessions from "Pending" to "Not Recorded" ---
# Each of these rows is the session dated 08/01/2026 for its group; as that
# date has now passed without a recording, it moves from the "upcoming"
# (style 6 / yellow "Pending") look to the "past due" (style 4 / pink
# "Not Recorded") look - matching the rows immediately above it that are
# already past their session date.
$rolloverRows = @(178, 205, 232, 259, 286, 313)
foreach ($row in $rolloverRows) {
    $srcFormatRow = $row - 3
    $src = $ws.Range("A" + $srcFormatRow + ":I" + $srcFormatRow)
    $dst = $ws.Range("A" + $row + ":I" + $row)
    $src.Copy()
    $dst.PasteSpecial(-4122)
    $ws.Cells.Item($row, 9).Value = "Not Recorded"
}

$excel.CutCopyMode = 0
